$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a text value into a cell while preventing Excel's
# automatic date/number inference from kicking in, and without leaving
# a residual cell style behind.
function Set-TextValue($rng, $value) {
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

# --- Edit existing rows -------------------------------------------------
# Row 280: fecha de pago corrected
Set-TextValue $ws.Range("B280") "24/12/2024"

# Row 292: fecha de pago corrected
Set-TextValue $ws.Range("B292") "10/01/2025"

# --- Insert three new rows before the current row 293 -------------------
# (this shifts the old rows 293-295 down to 296-298, matching the diff)
$ws.Rows(293).Insert()
$ws.Rows(293).Insert()
$ws.Rows(293).Insert()

# New row 293
Set-TextValue $ws.Range("B293") "24/12/2024"
Set-TextValue $ws.Range("D293") "Descuento - Anticipo"
Set-TextValue $ws.Range("E293") "Ivonne Mancipe"
Set-TextValue $ws.Range("G293") "Descuento"
$ws.Range("H293").Value = 0
$ws.Range("I293").Value = 0
$ws.Range("J293").Value = 0
$ws.Range("K293").Value = 0
$ws.Range("L293").Value = -20000

# New row 294
Set-TextValue $ws.Range("B294") "24/12/2024"
Set-TextValue $ws.Range("D294") "Descuento - Anticipo"
Set-TextValue $ws.Range("E294") "Ivonne Mancipe"
Set-TextValue $ws.Range("G294") "Descuento"
$ws.Range("H294").Value = 0
$ws.Range("I294").Value = 0
$ws.Range("J294").Value = 0
$ws.Range("K294").Value = 0
$ws.Range("L294").Value = -200000

# New row 295
Set-TextValue $ws.Range("B295") "24/12/2024"
Set-TextValue $ws.Range("D295") "Descuento - 4 Aluerzos"
Set-TextValue $ws.Range("E295") "Ivonne Mancipe"
Set-TextValue $ws.Range("G295") "Descuento"
$ws.Range("H295").Value = 0
$ws.Range("I295").Value = 0
$ws.Range("J295").Value = 0
$ws.Range("K295").Value = 0
$ws.Range("L295").Value = -60000
